# Weekly fruit/vegetable data update:
# Insert a new daily record as row 163 (pushing the existing rows 163-200
# down to 164-201) and populate it with the newest reading.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 163..200 down by one row so a new row 163 can be inserted
# above them (xlShiftDown = -4121).
$ws.Rows.Item(163).Insert(-4121)

# Populate the newly inserted row 163 with the new data point.
$ws.Cells.Item(163, 1).Value  = 9
$ws.Cells.Item(163, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(163, 3).Value  = "Metropolitana"
$ws.Cells.Item(163, 4).Value  = 45258
$ws.Cells.Item(163, 5).Value  = 13
$ws.Cells.Item(163, 6).Value  = 100112022
$ws.Cells.Item(163, 7).Value  = "Arveja Verde"
$ws.Cells.Item(163, 8).Value  = "Sin especificar"
$ws.Cells.Item(163, 9).Value  = "Primera"
$ws.Cells.Item(163, 10).Value = 40
$ws.Cells.Item(163, 11).Value = 21000
$ws.Cells.Item(163, 12).Value = 21000
$ws.Cells.Item(163, 13).Value = 21000
$ws.Cells.Item(163, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(163, 15).Value = "Región del Maule"
$ws.Cells.Item(163, 16).Value = 840
$ws.Cells.Item(163, 17).Value = 25
$ws.Cells.Item(163, 18).Value = "Hortaliza"
